$wb = $excel.ActiveWorkbook

# Color used by the workbook's "HyperLink" style (RGB FF6495ED), expressed
# as the BGR integer that Excel's Font.Color expects.
$hyperlinkColor = 15570276

function Set-HandbackRow($ws, $row, $statusText, $srcDisplay, $srcUrl, $xlfDisplay, $xlfUrl, $handbackDateTime) {
    # Column B: Status -> now handed back and in sync with en-US.
    $ws.Range("B$row").Value = $statusText

    # Column E: Latest Target File (new hyperlink, mirrors column A).
    $eCell = $ws.Range("E$row")
    $eCell.Value = $srcDisplay
    $ws.Hyperlinks.Add($eCell, $srcUrl, "", "", $srcDisplay) | Out-Null
    $eCell.Font.Underline = $true
    $eCell.Font.Color = $hyperlinkColor

    # Column F: Latest Handback File (new hyperlink, mirrors column C).
    $fCell = $ws.Range("F$row")
    $fCell.Value = $xlfDisplay
    $ws.Hyperlinks.Add($fCell, $xlfUrl, "", "", $xlfDisplay) | Out-Null
    $fCell.Font.Underline = $true
    $fCell.Font.Color = $hyperlinkColor

    # Column G: Latest Handback DateTime -> populate with the real timestamp.
    $ws.Range("G$row").Value = $handbackDateTime
}

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

Set-HandbackRow $wsZh 2 "Handed back: in sync with en-US" `
    "b3a76179-4237-41de-83e8-5e751a22958c.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/8753696ba97ccff6c58d455f230fa71c3ef1d7da/e2e/b3a76179-4237-41de-83e8-5e751a22958c.md" `
    "b3a76179-4237-41de-83e8-5e751a22958c.e3e9918caacf804b2dc7522769fa9ef66faabcd1.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/544e322317031b1c01f7facf71ea8362568b547f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/b3a76179-4237-41de-83e8-5e751a22958c.e3e9918caacf804b2dc7522769fa9ef66faabcd1.zh-cn.xlf" `
    "2016-03-02 15:37:48"

Set-HandbackRow $wsZh 3 "Handed back: in sync with en-US" `
    "ea31f9a5-27b9-4331-94c3-1c02177394ae.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/8753696ba97ccff6c58d455f230fa71c3ef1d7da/e2e/ea31f9a5-27b9-4331-94c3-1c02177394ae.md" `
    "ea31f9a5-27b9-4331-94c3-1c02177394ae.952e58fdebdd0751994f46d438059cd1e221ee9d.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/544e322317031b1c01f7facf71ea8362568b547f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ea31f9a5-27b9-4331-94c3-1c02177394ae.952e58fdebdd0751994f46d438059cd1e221ee9d.zh-cn.xlf" `
    "2016-03-02 15:37:48"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

Set-HandbackRow $wsDe 2 "Handed back: in sync with en-US" `
    "b3a76179-4237-41de-83e8-5e751a22958c.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/8753696ba97ccff6c58d455f230fa71c3ef1d7da/e2e/b3a76179-4237-41de-83e8-5e751a22958c.md" `
    "b3a76179-4237-41de-83e8-5e751a22958c.e3e9918caacf804b2dc7522769fa9ef66faabcd1.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/05c4e77b5495e85a17b10e64b74593087ccaf617/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/b3a76179-4237-41de-83e8-5e751a22958c.e3e9918caacf804b2dc7522769fa9ef66faabcd1.de-de.xlf" `
    "2016-03-02 15:38:07"

Set-HandbackRow $wsDe 3 "Handed back: in sync with en-US" `
    "ea31f9a5-27b9-4331-94c3-1c02177394ae.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/8753696ba97ccff6c58d455f230fa71c3ef1d7da/e2e/ea31f9a5-27b9-4331-94c3-1c02177394ae.md" `
    "ea31f9a5-27b9-4331-94c3-1c02177394ae.952e58fdebdd0751994f46d438059cd1e221ee9d.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/05c4e77b5495e85a17b10e64b74593087ccaf617/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ea31f9a5-27b9-4331-94c3-1c02177394ae.952e58fdebdd0751994f46d438059cd1e221ee9d.de-de.xlf" `
    "2016-03-02 15:38:07"

Write-Output "Handback report generated"
